{"js": "// Office.js (Word JavaScript API) script.\n// This is the body of: async (context) => { ... }\n//\n// Change being applied (from the diff):\n//  1. Delete the bullet paragraph \"When pressing T, the player should be\n//     teleported to the tutorial room, where one can try out key mappings\n//     we have.\" entirely.\n//  2. Fill in the text of the (previously empty) trailing bullet paragraph\n//     and append six brand-new bullet paragraphs after it, each using the\n//     same numbered-list formatting (numId 1), right before the final\n//     blank, non-list paragraph that ends the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- Step 1: remove the obsolete \"pressing T / tutorial room\" bullet ---\nconst removedText =\n  \"When pressing T, the player should be teleported to the tutorial room, \" +\n  \"where one can try out key mappings we have.\";\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === removedText) {\n    targetIndex = i;\n    break;\n  }\n}\nif (targetIndex !== -1) {\n  paragraphs.items[targetIndex].delete();\n  await context.sync();\n}\n\n// --- Step 2: locate the empty bullet paragraph right before the final ---\n// --- blank paragraph, and use it + new paragraphs for the new bullets ---\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText =\n  \"Shooting reimus bullet into walls should show a breaking animation\";\n\nlet anchorIndex = -1;\nfor (let i = 0; i < refreshedParagraphs.items.length; i++) {\n  if (refreshedParagraphs.items[i].text === anchorText) {\n    anchorIndex = i;\n    break;\n  }\n}\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate anchor paragraph for insertion.\");\n}\n\n// The paragraph right after the anchor is the empty bulleted list item\n// that currently has no text.\nconst emptyBulletParagraph = refreshedParagraphs.items[anchorIndex + 1];\n\nconst newBulletTexts = [\n  \"When starting a new game, or seeing a boss, the player should be able to interact with the dialogue box. Press space or mouse left to go to the next dialogue.\",\n  \"When pressing escape, instead of directly going out of the game, the player should be able to see the pause menu, where they can restart, resume, or exit to main menu.\",\n  \"When in the main menu, the player should be able to resume a game they are playing, or start a new game, or quit.\",\n  \"When enemies drop items that can be bought with coins, if players have enough coins, they should be able to buy it with \\u201cE\\u201d.\",\n  \"When entering a newly explored room, fog of war should expand on the current Reimu location to the entire room.\",\n  \"When a player tries to enter a room, doors should be unlocked. When the player is physically in the room, the door should be locked until the player kills all enemies in that room.\",\n  \"When players enter the main menu, they should hear a particular sound track. When the player goes in the room, they should hear another track that's different from the main menu.\",\n];\n\n// Fill the existing empty bullet paragraph with the first sentence.\nemptyBulletParagraph.insertText(newBulletTexts[0], \"End\");\nawait context.sync();\n\n// Insert the remaining bullet paragraphs after it (in order), each one\n// inherits the numbered-list formatting from the paragraph it follows.\nlet previousParagraph = emptyBulletParagraph;\nfor (let i = 1; i < newBulletTexts.length; i++) {\n  const inserted = previousParagraph.insertParagraph(newBulletTexts[i], \"After\");\n  await context.sync();\n  previousParagraph = inserted;\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $d is the active document ($word.ActiveDocument).\n#\n# Change being applied (from the diff):\n#  1. Delete the bullet paragraph \"When pressing T, the player should be\n#     teleported to the tutorial room, where one can try out key mappings\n#     we have.\" entirely.\n#  2. Fill in the text of the (previously empty) trailing bullet paragraph\n#     and append six brand-new bullet paragraphs after it, each using the\n#     same numbered-list formatting (numId 1), right before the final\n#     blank, non-list paragraph that ends the document.\n\n$d = $word.ActiveDocument\n\n$removedText = \"When pressing T, the player should be teleported to the tutorial room, where one can try out key mappings we have.\"\n$anchorText = \"Shooting reimus bullet into walls should show a breaking animation\"\n\n# --- Step 1: remove the obsolete \"pressing T / tutorial room\" bullet ---\n$toDelete = $null\nforeach ($p in $d.Paragraphs) {\n    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($ptext -eq $removedText) {\n        $toDelete = $p\n        break\n    }\n}\nif ($toDelete -ne $null) {\n    $toDelete.Range.Delete()\n}\n\n# --- Step 2: re-locate the anchor paragraph (indices shifted after the ---\n# --- delete above), then find the empty bullet paragraph right after it ---\n$anchorIndex = -1\n$idx = 1\nforeach ($p in $d.Paragraphs) {\n    $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($ptext -eq $anchorText) {\n        $anchorIndex = $idx\n        break\n    }\n    $idx++\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate anchor paragraph for insertion.\"\n}\n\n$newBulletTexts = @(\n    \"When starting a new game, or seeing a boss, the player should be able to interact with the dialogue box. Press space or mouse left to go to the next dialogue.\",\n    \"When pressing escape, instead of directly going out of the game, the player should be able to see the pause menu, where they can restart, resume, or exit to main menu.\",\n    \"When in the main menu, the player should be able to resume a game they are playing, or start a new game, or quit.\",\n    (\"When enemies drop items that can be bought with coins, if players have enough coins, they should be able to buy it with \" + [char]0x201C + \"E\" + [char]0x201D + \".\"),\n    \"When entering a newly explored room, fog of war should expand on the current Reimu location to the entire room.\",\n    \"When a player tries to enter a room, doors should be unlocked. When the player is physically in the room, the door should be locked until the player kills all enemies in that room.\",\n    \"When players enter the main menu, they should hear a particular sound track. When the player goes in the room, they should hear another track that's different from the main menu.\"\n)\n\n# The paragraph right after the anchor is the empty bulleted list item\n# that currently has no text - reuse it for the first new sentence.\n$bulletIndex = $anchorIndex + 1\n$firstPara = $d.Paragraphs.Item($bulletIndex)\n$firstPara.Range.Text = $newBulletTexts[0]\n\n# Insert the remaining bullet paragraphs after it (in order). Each fresh\n# InsertParagraphAfter() inherits the numbered-list formatting (numId 1)\n# from the paragraph it follows. Re-fetch paragraphs by index (rather than\n# chaining .Next()) so each reference stays valid after the mutation.\nfor ($i = 1; $i -lt $newBulletTexts.Length; $i++) {\n    $srcIndex = $bulletIndex + $i - 1\n    $srcPara = $d.Paragraphs.Item($srcIndex)\n    $srcPara.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($srcIndex + 1)\n    $newPara.Range.Text = $newBulletTexts[$i]\n}\n"}
